$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 6 so the old row 6 (quotation) shifts down to row 7.
$ws.Rows.Item(6).Insert()

# New row 6: a refund record linked to base_id 20251116-002
$ws.Cells.Item(6, 1).Value = "20251116-002"

# Date column would otherwise get auto-parsed into a date serial; force text.
$ws.Cells.Item(6, 2).NumberFormat = "@"
$ws.Cells.Item(6, 2).Value = "2025-11-16"
$ws.Cells.Item(6, 2).Style = "Normal"

$ws.Cells.Item(6, 3).Value = "r"
$ws.Cells.Item(6, 4).Value = "R-20251116-20251116-002-1"
$ws.Cells.Item(6, 5).Value = 0
$ws.Cells.Item(6, 6).Value = "Fahad Ahmed Mohammed"
$ws.Cells.Item(6, 7).Value = ""
$ws.Cells.Item(6, 8).Value = "Abu Dhabi - Al Shamkha"
$ws.Cells.Item(6, 9).Value = ""

# Update the shifted-down row 7 (was row 6): base_id + amount changed, contact info cleared
$ws.Cells.Item(7, 1).Value = "20251116-007"
$ws.Cells.Item(7, 5).Value = 4680
$ws.Cells.Item(7, 6).Value = ""
$ws.Cells.Item(7, 7).Value = ""
